# 10Th - MB for single stock and added new group
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (old B:E -> shifts to E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# New date headers for the newly inserted columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns with "UN" for all existing data rows
$ws.Range("B2:D27").Value = "UN"

# Add two new analyst rows at the bottom of the table
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
